# Update the "想去人数" (interested-count) figures in column F on the
# "展览" and "全部类型" sheets to match the refreshed scrape.

$wb = $excel.ActiveWorkbook

# row => new value, applied on sheet "展览"
$updatesExhibition = @{
    3  = 16465
    4  = 22
    5  = 740
    6  = 15626
    8  = 9288
    9  = 502
    13 = 229
    18 = 629
    22 = 1160
    26 = 540
    28 = 46
    37 = 5735
    38 = 5255
}

# row => new value, applied on sheet "全部类型" (same events, shifted rows)
$updatesAllTypes = @{
    3  = 16465
    4  = 22
    5  = 740
    6  = 15626
    8  = 9288
    9  = 502
    13 = 229
    18 = 629
    22 = 1160
    26 = 540
    28 = 46
    39 = 5735
    41 = 5255
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $updatesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAllTypes.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $updatesAllTypes[$row]
}
